$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 186 (existing rows 186-264 shift down to 187-265).
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new price-record data.
$ws.Range("A186").Value = 7
$ws.Range("B186").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C186").Value = "Ñuble"
$ws.Range("D186").Value = 45134
$ws.Range("E186").Value = 16
$ws.Range("F186").Value = "Fruta"
$ws.Range("G186").Value = 100109
$ws.Range("H186").Value = "Uva"
$ws.Range("I186").Value = 100109001
$ws.Range("J186").Value = "Uva"
$ws.Range("K186").Value = "Crimpson Seedless"
$ws.Range("L186").Value = "Primera"
$ws.Range("M186").Value = 50
$ws.Range("N186").Value = 12000
$ws.Range("O186").Value = 12000
$ws.Range("P186").Value = 12000
$ws.Range("Q186").Value = "$/bandeja 8 kilos"
$ws.Range("R186").Value = "Región de O'Higgins"
$ws.Range("S186").Value = 1500
$ws.Range("T186").Value = 8
